$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.934.11'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '2.211.47'
$ws.Range("E3").Value = '  -0.72%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '289.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.511'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.72%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0776'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.110'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.09%  '

$ws.Range("D14").Value = '2.551.97'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.31%  '

$ws.Range("D16").Value = '2.211.78'
$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.726'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '

$ws.Range("D18").Value = '39.866.33'
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.48%  '

$ws.Range("D20").Value = '0.0₃0882'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '155.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0716'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.83'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.111'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0984'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.82%  '

$ws.Range("D42").Value = '2.093.90'
$ws.Range("E42").Value = '  +7.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0268'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.98%  '

$ws.Range("D48").Value = '2.425.16'
$ws.Range("E48").Value = '  -0.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '68.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.65%  '
